$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last row (row 15), shifting the
# existing "Willian Massami Watanabe" row down to row 17.
$ws.Rows("15:16").Insert()

# Row 15: Reginaldo Fidelis
$ws.Range("A15").Value = "Reginaldo Fidelis"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0

# Row 16: Rogerio Tondato
$ws.Range("A16").Value = "Rogerio Tondato"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 6
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0

$wb.Save()
